# Fix two typos in the .budgt Ruby Terminal App deck.
#
# Slide 1, subtitle shape: "were just saving it for later." -> "were just
# saving up for lypo."
# Slide 2, body shape: "...family members ar christmas!" -> "...family
# members at christmas!"
#
# NOTE: TextRange.Text normalizes curly quotes to straight quotes when
# read back, so the search key below uses a straight quote; the
# replacement text uses the real typographic quotes (“ ”) so the saved
# OOXML keeps the original character.

$p = $ppt.ActivePresentation

# --- Slide 1 -------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$shp1 = $s1.Shapes.Item(2)
$tr1 = $shp1.TextFrame.TextRange

$oldText1 = "The missing `"e`" is not a typo, were just saving it for later."
$newText1 = "The missing “e” is not a typo, were just saving up for lypo."

$idx1 = $tr1.Text.IndexOf($oldText1)
if ($idx1 -lt 0) {
    throw "Slide 1 target text not found"
}
$run1 = $tr1.Characters($idx1 + 1, $oldText1.Length)
$run1.Text = $newText1

# --- Slide 2 -------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item(4)
$tr2 = $shp2.TextFrame.TextRange

$oldText2 = "This will help you save so you can spend all your cold hard savings on ungrateful family members ar christmas!"
$newText2 = "This will help you save so you can spend all your cold hard savings on ungrateful family members at christmas!"

$idx2 = $tr2.Text.IndexOf($oldText2)
if ($idx2 -lt 0) {
    throw "Slide 2 target text not found"
}
$run2 = $tr2.Characters($idx2 + 1, $oldText2.Length)
$run2.Text = $newText2
